# Marketing plan connection matrix - add new channels (Update trailer /
# Update landing page) as columns, and new outreach rows (Email to existing
# creators / Update trailer / Update landing page / Email to new youtubers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (K, L) on row 1 ---
$ws.Range("K1").Value = "Update trailer"
$ws.Range("L1").Value = "Update landing page"

# --- Row 2 ("In game") gains a couple of new cells / value changes ---
$ws.Range("C2").Value = "not needed"
$ws.Range("E2").Value = "mentioned"
$ws.Range("F2").Value = "mentioned"
$ws.Range("J2").Value = "not important enough"

# Excel records "mentioned" (no longer "could be mentioned") with an
# explicit (no-fill) style on E2/F2, where previously a fill highlighted
# "could be mentioned".
$ws.Range("E2").Interior.Pattern = -4142
$ws.Range("F2").Interior.Pattern = -4142

# --- New rows 13-16 ---
$ws.Range("A13").Value = "Email to existing creators"
$ws.Range("B13").Value = "not needed"
$ws.Range("C13").Value = "not needed"
$ws.Range("D13").Value = "not needed"
$ws.Range("E13").Value = "indirect via trailer"
$ws.Range("F13").Value = "indirect via trailer"
$ws.Range("G13").Value = "not needed"
$ws.Range("H13").Value = "not important enough"
$ws.Range("I13").Value = "not needed"
$ws.Range("J13").Value = "not needed"
$ws.Range("K13").Value = "mentioned"
$ws.Range("L13").Value = "indirect via trailer"

$ws.Range("A14").Value = "Update trailer"
$ws.Range("A15").Value = "Update landing page"
$ws.Range("A16").Value = "Email to new youtubers"

# --- Column widths (best-fit on the widened columns) ---
$ws.Columns.Item(1).ColumnWidth = 20.94
$ws.Columns.Item(11).ColumnWidth = 11.5
$ws.Columns.Item(12).ColumnWidth = 16.94

# --- Selection ends on F16, matching the author's last edit location ---
$ws.Range("F16").Select() | Out-Null
